$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the two new "BBR" rows: insert 3 blank rows starting at
#    row 38 (this pushes the "Line Graph:" label down to row 41, and the
#    whole second table down from rows 39-63 to rows 42-66).
# ---------------------------------------------------------------------------
$ws.Range("A38:A40").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the previously-zeroed TGSM / Graph2Vec rows (36 = mean,
#    37 = initial) with the real averagePrecision numbers.
# ---------------------------------------------------------------------------
$row36 = @(0.59, 0.45, 0.55, 0.28, 0.53, 0.5, 0.5, 0.55, 0.5)
for ($i = 0; $i -lt $row36.Length; $i++) {
    $ws.Cells.Item(36, 4 + $i).Value = $row36[$i]
}

$row37 = @(0.5, 0.47, 0.5, 0, 0.51, 0.43, 0.48, 0.5, 0.42)
for ($i = 0; $i -lt $row37.Length; $i++) {
    $ws.Cells.Item(37, 4 + $i).Value = $row37[$i]
}

# ---------------------------------------------------------------------------
# 3. New rows 38/39: BBR / Graph2Vec, mean + initial reference vectors.
# ---------------------------------------------------------------------------
$ws.Cells.Item(38, 1).Value = "BBR"
$ws.Cells.Item(38, 2).Value = "Graph2Vec"
$ws.Cells.Item(38, 3).Value = "mean"
$row38 = @(0.08, 0.05, 0.08, 0.04, 0.08, 0.08, 0.12, 0.08, 0.05)
for ($i = 0; $i -lt $row38.Length; $i++) {
    $ws.Cells.Item(38, 4 + $i).Value = $row38[$i]
}

$ws.Cells.Item(39, 1).Value = "BBR"
$ws.Cells.Item(39, 2).Value = "Graph2Vec"
$ws.Cells.Item(39, 3).Value = "initial"
$row39 = @(0.06, 0.07, 0.06, 0, 0.06, 0.07, 0.05, 0.06, 0.06)
for ($i = 0; $i -lt $row39.Length; $i++) {
    $ws.Cells.Item(39, 4 + $i).Value = $row39[$i]
}

# Row 40 must stay completely blank (it's the spacer row before the
# "Line Graph:" caption that got pushed down to row 41).
$ws.Rows.Item(40).Clear()

# ---------------------------------------------------------------------------
# 4. Grow Table1 so it covers the two new data rows (was A1:L38, now
#    A1:L39 to match the extra BBR rows that were inserted right before
#    its old bottom edge).
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:L39"))

# ---------------------------------------------------------------------------
# 5. Restore the view: scroll back to the top and move the selection to
#    where the author last left off.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C68").Select()
